# Tripadvisor New Orleans shard 58 update:
#   1. Insert a new "State" column into the hotel_info sheet (between
#      Hotel_Name and City) and populate it with "Louisiana".
#   2. Reorder the worksheets so review_info precedes hotel_info.

$wb = $excel.ActiveWorkbook

$wsHotel  = $wb.Worksheets.Item("hotel_info")
$wsReview = $wb.Worksheets.Item("review_info")

# --- 1. Add the "State" column -------------------------------------------
# hotel_info columns (before): A=STR, B=Hotel_Name, C=City, D=Zip,
#   E=TA_ReviewURL, F=Tripadvisor_Hotel_Name, G=English_Reviews_num,
#   H=Local_Rank, I=Total_Reviews_num
# Insert a blank column at C so State sits right after Hotel_Name.
$wsHotel.Columns.Item(3).Insert()
$wsHotel.Cells.Item(1, 3).Value = "State"
$wsHotel.Cells.Item(2, 3).Value = "Louisiana"

# --- 2. Reorder sheets: review_info, then hotel_info ----------------------
$wsHotel.Move($null, $wsReview)
